# Adding new "Categories" Tab  (#269)
# Mirrors the header formatting of the "ChallengeQuestion" sheet (the
# previously active tab) - same row-1 key-colour headers (s=12/13/4/5)
# and row-2 description styles (s=7/8).

$wb = $excel.ActiveWorkbook

# Touch the old active sheet's header block first (A1:H2) - this is the
# range the new sheet's formatting/content is modelled on.
$fmtSrc = $wb.Worksheets.Item("ChallengeQuestion")
$fmtSrc.Range("A1:H2").Select() | Out-Null

# ---------------------------------------------------------------------
# 1. Create the new worksheet at the very end of the workbook and name it.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Categories"

# ---------------------------------------------------------------------
# 2. Column widths (only the columns that deviate from the sheet default).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(5).ColumnWidth = 17.666666666666668
$ws.Columns.Item(6).ColumnWidth = 11.5
$ws.Columns.Item(7).ColumnWidth = 16.833333333333332

# ---------------------------------------------------------------------
# 3. Row heights for the header rows.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 18
$ws.Rows.Item(2).RowHeight = 113

# ---------------------------------------------------------------------
# 4. Row 1 - title + common column-colour-key headers.
# ---------------------------------------------------------------------
$fmtSrc.Range("A1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$ws.Range("A1").Value = "Categories"

$fmtSrc.Range("C1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").Value = "PrimaryKeyInRed"

$fmtSrc.Range("D1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = "PrimaryAndForeignKey Orange"

$fmtSrc.Range("E1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Value = "ForeignKey Brown"

# ---------------------------------------------------------------------
# 5. Row 2 - column descriptions.
# ---------------------------------------------------------------------
$fmtSrc.Range("A2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Value = "Start date from which the data will be valid"

$fmtSrc.Range("B2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Value = "End date until which the data will be valid"

$fmtSrc.Range("B2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Value = "The identifier which defines the CaseType code.`nMaxLength: 70"

$fmtSrc.Range("B2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Value = "Unique ID that identifies the category`nMaxLength: 70"

$fmtSrc.Range("B2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Value = "Label that gets displayed in the UI. MaxLength: 70"

$fmtSrc.Range("B2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Value = "Specifies the display order for the category."

$fmtSrc.Range("B2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Value = "Identifies the parent category id for a sub-category."

$fmtSrc.Range("B2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 6. Row 3 - actual field/column names (plain, unstyled cells).
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "LiveFrom"
$ws.Range("B3").Value = "LiveTo"
$ws.Range("C3").Value = "CaseTypeID"
$ws.Range("D3").Value = "CategoryID"
$ws.Range("E3").Value = "CategoryLabel"
$ws.Range("F3").Value = "DisplayOrder"
$ws.Range("G3").Value = "ParentCategoryID"

# ---------------------------------------------------------------------
# 7. Leave the new tab as the active / selected sheet & cell.
# ---------------------------------------------------------------------
$ws.Range("O23").Select() | Out-Null
